$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.299.18"
$ws.Range("E2").Value = "  -0.87%  "

$ws.Range("D3").Value = "3.093.99"
$ws.Range("E3").Value = "  -0.28%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "524.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.24%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "3.093.71"
$ws.Range("E8").Value = "  -0.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.449"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.01%  "

$ws.Range("E11").Value = "  -1.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.394"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.43%  "

$ws.Range("D13").Value = "3.623.54"
$ws.Range("E13").Value = "  -0.28%  "

$ws.Range("E14").Value = "  +2.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.23"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.12%  "

$ws.Range("E16").Value = "  -1.24%  "

$ws.Range("D17").Value = "57.321.53"
$ws.Range("E17").Value = "  -0.90%  "

$ws.Range("D18").Value = "3.093.98"
$ws.Range("E18").Value = "  -0.16%  "

$ws.Range("E19").Value = "  -3.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.61%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.61%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "346.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.69%  "

$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.83%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.498"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.86%  "

$ws.Range("E26").Value = "  -2.50%  "

$ws.Range("E27").Value = "  +0.03%  "

$ws.Range("D28").Value = "0.0₃0885"
$ws.Range("E28").Value = "  -3.98%  "

$ws.Range("E29").Value = "  -0.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.87%  "

$ws.Range("E31").Value = "  -0.47%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.58%  "

$ws.Range("E33").Value = "  -1.67%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.06%  "

$ws.Range("E35").Value = "  -4.12%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "158.70"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.60%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.74%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.66"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.43%  "

$ws.Range("E39").Value = "  -2.39%  "

$ws.Range("E40").Value = "  +5.45%  "

$ws.Range("E41").Value = "  -1.23%  "

$ws.Range("E42").Value = "  +2.58%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.695"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.97%  "

$ws.Range("D44").Value = "2.367.67"
$ws.Range("E44").Value = "  +3.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.00%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0265"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.90%  "

$ws.Range("E48").Value = "  -1.81%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.45%  "

$ws.Range("E50").Value = "  -4.62%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.754"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.60%  "
